# Column C ("Förändrad") holds a "last changed" date serial value that was
# bumped by one day (46061 -> 46062, i.e. 2026-02-08 -> 2026-02-09) for every
# data row (rows 2 through 431) in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C431").Value = 46062
